$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to text,
# otherwise Excel auto-converts the assignment to a numeric value (matching
# genuine Excel type-inference behavior) and the original formatting (e.g.
# trailing zeros) would be lost. We briefly mark the cell as Text, assign the
# value, then restore the default "Normal" style so no stray formatting remains.

$ws.Range("D2").Value = "26.644.48"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.598.11"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.822.60"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.601.78"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "26.629.72"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "1.275.17"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  -7.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  +18.28%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "1.735.19"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "
